$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value (45175 -> 2023-09-06)
# that needs to be updated to 45177 (-> 2023-09-08) for every data row
# (rows 2 through 301).
$lastRow = 301
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45177
